$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-11 from 45183 to 45184
$ws.Range("C2:C11").Value = 45184
